$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Trim trailing/leading whitespace off the "Projekt-XXX" labels in column A
# (rows 10-26 hold the per-project task rows; A25/A26 are already clean).
for ($r = 10; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2.Trim()
}

# Restore the selection the workbook was saved with.
[void]$ws.Range("G23").Select()
